# "Generate Report for Handback"
# The localization-status report is regenerated after a handback: the
# "Ready for handoff" status becomes "Handed back: in sync with en-US",
# and the per-language sheets get the freshly generated handback file
# name + timestamp recorded for each source file (a.md / b.md).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/000222a3b3e70f93181ffa723c3cbd634e306f77/e2e"

# ---------------------------------------------------------------------
# Overview sheet: status columns for zh-cn (E) and de-de (F) move from
# "Ready for handoff" to "Handed back: in sync with en-US", and those
# two columns widen to fit the longer text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: Status column, the Latest Target File / Latest Handback
# File / Latest Handback DateTime columns get filled in for both rows.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-05 16:44:18"
$wsZh.Range("K3").Value = "2016-09-05 16:44:18"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# Rebuild the hyperlinks in row order (A2, I2, A3, I3) so the new
# "Latest Target File" links (column I) line up with the rest.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$baseUrl/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$baseUrl/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$baseUrl/b.md", [System.Type]::Missing, [System.Type]::Missing, "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$baseUrl/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md")

# ---------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, but with its own handback
# file name + its own (later) handback timestamp.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-05 16:44:25"
$wsDe.Range("K3").Value = "2016-09-05 16:44:25"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$baseUrl/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$baseUrl/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$baseUrl/b.md", [System.Type]::Missing, [System.Type]::Missing, "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$baseUrl/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md")
